$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates: volume/issue number and reporting week date range
$ws.Range("A8").Value = "Volume 31   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/22/2024  Through  4/28/2024"

# Crime statistics data updates (rows 14-30)
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("I14").Value = 5
$ws.Range("K14").Value = 150
$ws.Range("L14").Value = 150
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -80
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = 50
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -85.714285714285
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -52.380952380952
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 61
$ws.Range("K16").Value = -18.032786885245
$ws.Range("L16").Value = -27.536231884058
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -92.732558139534
$ws.Range("C17").Value = 6
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = -41.463414634146
$ws.Range("I17").Value = 96
$ws.Range("J17").Value = 131
$ws.Range("K17").Value = -26.717557251908
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 12.941176470588
$ws.Range("N17").Value = -72.492836676217
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -65
$ws.Range("I18").Value = 42
$ws.Range("J18").Value = 48
$ws.Range("K18").Value = -12.5
$ws.Range("L18").Value = -31.147540983606
$ws.Range("M18").Value = -51.724137931034
$ws.Range("N18").Value = -91.764705882352
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -26.470588235294
$ws.Range("I19").Value = 119
$ws.Range("J19").Value = 136
$ws.Range("K19").Value = -12.5
$ws.Range("L19").Value = -13.138686131386
$ws.Range("M19").Value = 52.564102564102
$ws.Range("N19").Value = -36.702127659574
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -80
$ws.Range("I20").Value = 24
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = -52
$ws.Range("L20").Value = -35.135135135135
$ws.Range("M20").Value = -35.135135135135
$ws.Range("N20").Value = -91.608391608391
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -43.75
$ws.Range("G21").Value = 129
$ws.Range("H21").Value = -44.186046511627
$ws.Range("I21").Value = 345
$ws.Range("J21").Value = 437
$ws.Range("K21").Value = -21.052631578947
$ws.Range("L21").Value = -15.647921760391
$ws.Range("M21").Value = -7.506702412868
$ws.Range("N21").Value = -83.284883720930
$ws.Range("G22").Value = 3
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = -50
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -53.846153846153
$ws.Range("I23").Value = 29
$ws.Range("J23").Value = 48
$ws.Range("K23").Value = -39.583333333333
$ws.Range("L23").Value = -6.451612903225
$ws.Range("M23").Value = 52.631578947368
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -46.666666666666
$ws.Range("F24").Value = 41
$ws.Range("G24").Value = 59
$ws.Range("H24").Value = -30.508474576271
$ws.Range("I24").Value = 221
$ws.Range("J24").Value = 235
$ws.Range("K24").Value = -5.957446808510
$ws.Range("L24").Value = -19.047619047619
$ws.Range("M24").Value = -15
$ws.Range("C25").Value = 1
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 7
$ws.Range("H25").Value = 40
$ws.Range("J25").Value = 24
$ws.Range("K25").Value = 62.5
$ws.Range("L25").Value = 50
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -30
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = -9.375
$ws.Range("I26").Value = 135
$ws.Range("J26").Value = 189
$ws.Range("K26").Value = -28.571428571428
$ws.Range("L26").Value = -16.149068322981
$ws.Range("M26").Value = -44.214876033057
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = 60
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "***.*"
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 15
$ws.Range("K28").Value = -31.818181818181
$ws.Range("L28").Value = -54.545454545454
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -66.666666666666
$ws.Range("I29").Value = 4
$ws.Range("J29").Value = 6
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = -20
$ws.Range("M29").Value = -81.818181818181
$ws.Range("N29").Value = -95.121951219512
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -66.666666666666
$ws.Range("I30").Value = 4
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -33.333333333333
$ws.Range("M30").Value = -76.470588235294
$ws.Range("N30").Value = -94.805194805194
